# Auto-generated update of cached market-price / profit columns (H:N)
# across multiple worksheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 2497.5
$ws.Cells.Item(76, 9).Value = 2497.5
$ws.Cells.Item(76, 11).Value = 2497.5
$ws.Cells.Item(76, 13).Value = -2182.5
$ws.Cells.Item(79, 8).Value = 2497.5
$ws.Cells.Item(79, 9).Value = 2497.5
$ws.Cells.Item(79, 11).Value = 2497.5
$ws.Cells.Item(79, 13).Value = -1405.5
$ws.Cells.Item(116, 8).Value = 873748.5
$ws.Cells.Item(116, 9).Value = 1164015.5
$ws.Cells.Item(116, 11).Value = 1164015.5
$ws.Cells.Item(116, 13).Value = -1160573.5
$ws.Cells.Item(132, 8).Value = 4298.6343
$ws.Cells.Item(132, 9).Value = 4298.6343
$ws.Cells.Item(132, 11).Value = 12895.9029
$ws.Cells.Item(132, 13).Value = -10365.9029
$ws.Cells.Item(137, 8).Value = 6191.0454
$ws.Cells.Item(137, 9).Value = 1727.0667
$ws.Cells.Item(137, 10).Value = 15756.714
$ws.Cells.Item(137, 11).Value = 5181.2001
$ws.Cells.Item(137, 12).Value = 47270.142
$ws.Cells.Item(137, 13).Value = -2631.2001
$ws.Cells.Item(137, 14).Value = -52370.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5860.2407
$ws.Cells.Item(32, 9).Value = 4787.8887
$ws.Cells.Item(32, 10).Value = 11222
$ws.Cells.Item(32, 11).Value = 4787.8887
$ws.Cells.Item(32, 12).Value = 11222
$ws.Cells.Item(32, 13).Value = -4500.8887
$ws.Cells.Item(32, 14).Value = -11796
$ws.Cells.Item(43, 8).Value = 47496.75
$ws.Cells.Item(43, 10).Value = 49995.668
$ws.Cells.Item(43, 12).Value = 49995.668
$ws.Cells.Item(43, 14).Value = -50621.668
$ws.Cells.Item(45, 8).Value = 92038
$ws.Cells.Item(45, 9).Value = 146733.33
$ws.Cells.Item(45, 11).Value = 146733.33
$ws.Cells.Item(45, 13).Value = -146356.33
$ws.Cells.Item(61, 8).Value = 4647.05
$ws.Cells.Item(61, 9).Value = 2680.7778
$ws.Cells.Item(61, 11).Value = 2680.7778
$ws.Cells.Item(61, 13).Value = -2468.7778
$ws.Cells.Item(74, 8).Value = 285072.1
$ws.Cells.Item(74, 9).Value = 695814.75
$ws.Cells.Item(74, 11).Value = 695814.75
$ws.Cells.Item(74, 13).Value = -694940.75
$ws.Cells.Item(77, 8).Value = 285072.1
$ws.Cells.Item(77, 9).Value = 695814.75
$ws.Cells.Item(77, 11).Value = 3479073.75
$ws.Cells.Item(77, 13).Value = -3474705.75
$ws.Cells.Item(110, 8).Value = 7847.4546
$ws.Cells.Item(110, 9).Value = 7369.8887
$ws.Cells.Item(110, 10).Value = 9996.5
$ws.Cells.Item(110, 11).Value = 7369.8887
$ws.Cells.Item(110, 12).Value = 9996.5
$ws.Cells.Item(110, 13).Value = -5324.8887
$ws.Cells.Item(110, 14).Value = -14086.5
$ws.Cells.Item(132, 8).Value = 3392.8462
$ws.Cells.Item(132, 9).Value = 2799.5
$ws.Cells.Item(132, 10).Value = 3656.5557
$ws.Cells.Item(132, 11).Value = 8398.5
$ws.Cells.Item(132, 12).Value = 10969.6671
$ws.Cells.Item(132, 13).Value = -5868.5
$ws.Cells.Item(132, 14).Value = -16029.6671
$ws.Cells.Item(136, 8).Value = 4647.05
$ws.Cells.Item(136, 9).Value = 2680.7778
$ws.Cells.Item(136, 11).Value = 8042.3334
$ws.Cells.Item(136, 13).Value = -5492.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 23813764
$ws.Cells.Item(20, 9).Value = 28740116
$ws.Cells.Item(20, 10).Value = 3065.5
$ws.Cells.Item(20, 11).Value = 28740116
$ws.Cells.Item(20, 12).Value = 3065.5
$ws.Cells.Item(20, 13).Value = -28739869
$ws.Cells.Item(20, 14).Value = -3559.5
$ws.Cells.Item(55, 8).Value = 49999
$ws.Cells.Item(55, 10).Value = 49999
$ws.Cells.Item(55, 12).Value = 49999
$ws.Cells.Item(55, 14).Value = -50545
$ws.Cells.Item(86, 8).Value = 4023.8
$ws.Cells.Item(86, 9).Value = 3851.7334
$ws.Cells.Item(86, 10).Value = 4540
$ws.Cells.Item(86, 11).Value = 3851.7334
$ws.Cells.Item(86, 12).Value = 4540
$ws.Cells.Item(86, 13).Value = -2728.7334
$ws.Cells.Item(86, 14).Value = -6786
$ws.Cells.Item(89, 8).Value = 4023.8
$ws.Cells.Item(89, 9).Value = 3851.7334
$ws.Cells.Item(89, 10).Value = 4540
$ws.Cells.Item(89, 11).Value = 19258.667
$ws.Cells.Item(89, 12).Value = 22700
$ws.Cells.Item(89, 13).Value = -13642.667
$ws.Cells.Item(89, 14).Value = -33932
$ws.Cells.Item(99, 8).Value = 204502
$ws.Cells.Item(99, 9).Value = 501755
$ws.Cells.Item(99, 11).Value = 501755
$ws.Cells.Item(99, 13).Value = -500257
$ws.Cells.Item(134, 8).Value = 4165.231
$ws.Cells.Item(134, 9).Value = 4226.8
$ws.Cells.Item(134, 11).Value = 12680.4
$ws.Cells.Item(134, 13).Value = -10145.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2498
$ws.Cells.Item(16, 9).Value = 2498
$ws.Cells.Item(16, 11).Value = 2498
$ws.Cells.Item(16, 13).Value = -2211
$ws.Cells.Item(31, 8).Value = 3473.7637
$ws.Cells.Item(31, 9).Value = 2852.0312
$ws.Cells.Item(31, 11).Value = 2852.0312
$ws.Cells.Item(31, 13).Value = -2557.0312
$ws.Cells.Item(34, 8).Value = 3473.7637
$ws.Cells.Item(34, 9).Value = 2852.0312
$ws.Cells.Item(34, 11).Value = 2852.0312
$ws.Cells.Item(34, 13).Value = -2650.0312
$ws.Cells.Item(58, 8).Value = 3429.0454
$ws.Cells.Item(58, 9).Value = 3916.3333
$ws.Cells.Item(58, 10).Value = 3246.3125
$ws.Cells.Item(58, 11).Value = 3916.3333
$ws.Cells.Item(58, 12).Value = 3246.3125
$ws.Cells.Item(58, 13).Value = -3713.3333
$ws.Cells.Item(58, 14).Value = -3652.3125
$ws.Cells.Item(107, 8).Value = 2500626.8
$ws.Cells.Item(107, 9).Value = 3125624
$ws.Cells.Item(107, 11).Value = 3125624
$ws.Cells.Item(107, 13).Value = -3123704
$ws.Cells.Item(113, 8).Value = 2498
$ws.Cells.Item(113, 9).Value = 2498
$ws.Cells.Item(113, 11).Value = 2498
$ws.Cells.Item(113, 13).Value = -328
$ws.Cells.Item(132, 8).Value = 12503176
$ws.Cells.Item(132, 9).Value = 15627845
$ws.Cells.Item(132, 10).Value = 4499.125
$ws.Cells.Item(132, 11).Value = 46883535
$ws.Cells.Item(132, 12).Value = 13497.375
$ws.Cells.Item(132, 13).Value = -46881005
$ws.Cells.Item(132, 14).Value = -18557.375
$ws.Cells.Item(134, 8).Value = 5013.6
$ws.Cells.Item(134, 9).Value = 5049.6
$ws.Cells.Item(134, 11).Value = 15148.8
$ws.Cells.Item(134, 13).Value = -12613.8
$ws.Cells.Item(136, 8).Value = 3429.0454
$ws.Cells.Item(136, 9).Value = 3916.3333
$ws.Cells.Item(136, 10).Value = 3246.3125
$ws.Cells.Item(136, 11).Value = 11748.9999
$ws.Cells.Item(136, 12).Value = 9738.9375
$ws.Cells.Item(136, 13).Value = -9198.999899999999
$ws.Cells.Item(136, 14).Value = -14838.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15216867
$ws.Cells.Item(70, 9).Value = 21830418
$ws.Cells.Item(70, 10).Value = 5700.1
$ws.Cells.Item(70, 11).Value = 21830418
$ws.Cells.Item(70, 12).Value = 5700.1
$ws.Cells.Item(70, 13).Value = -21830148
$ws.Cells.Item(70, 14).Value = -6240.1
$ws.Cells.Item(73, 8).Value = 15216867
$ws.Cells.Item(73, 9).Value = 21830418
$ws.Cells.Item(73, 10).Value = 5700.1
$ws.Cells.Item(73, 11).Value = 21830418
$ws.Cells.Item(73, 12).Value = 5700.1
$ws.Cells.Item(73, 13).Value = -21829482
$ws.Cells.Item(73, 14).Value = -7572.1
$ws.Cells.Item(123, 8).Value = 57498.668
$ws.Cells.Item(123, 9).Value = 40000
$ws.Cells.Item(123, 10).Value = 74997.336
$ws.Cells.Item(123, 11).Value = 40000
$ws.Cells.Item(123, 12).Value = 74997.336
$ws.Cells.Item(123, 13).Value = -37550
$ws.Cells.Item(123, 14).Value = -79897.336
$ws.Cells.Item(132, 8).Value = 3257.1482
$ws.Cells.Item(132, 9).Value = 3472.3
$ws.Cells.Item(132, 10).Value = 2642.4285
$ws.Cells.Item(132, 11).Value = 10416.9
$ws.Cells.Item(132, 12).Value = 7927.2855
$ws.Cells.Item(132, 13).Value = -7886.900000000001
$ws.Cells.Item(132, 14).Value = -12987.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1999.75
$ws.Cells.Item(46, 10).Value = 1800.2
$ws.Cells.Item(46, 12).Value = 1800.2
$ws.Cells.Item(46, 14).Value = -2176.2
$ws.Cells.Item(82, 8).Value = 1739.5
$ws.Cells.Item(82, 10).Value = 494
$ws.Cells.Item(82, 12).Value = 494
$ws.Cells.Item(82, 14).Value = -1216
$ws.Cells.Item(85, 8).Value = 1739.5
$ws.Cells.Item(85, 10).Value = 494
$ws.Cells.Item(85, 12).Value = 494
$ws.Cells.Item(85, 14).Value = -2990
$ws.Cells.Item(122, 8).Value = 3178.4348
$ws.Cells.Item(122, 9).Value = 3341.7058
$ws.Cells.Item(122, 10).Value = 2715.8333
$ws.Cells.Item(122, 11).Value = 10025.1174
$ws.Cells.Item(122, 12).Value = 8147.499899999999
$ws.Cells.Item(122, 13).Value = -7575.117400000001
$ws.Cells.Item(122, 14).Value = -13047.4999
$ws.Cells.Item(132, 8).Value = 4409.1177
$ws.Cells.Item(132, 9).Value = 3588.122
$ws.Cells.Item(132, 11).Value = 10764.366
$ws.Cells.Item(132, 13).Value = -8234.366
$ws.Cells.Item(136, 8).Value = 5303.7646
$ws.Cells.Item(136, 9).Value = 6264.222
$ws.Cells.Item(136, 10).Value = 4223.25
$ws.Cells.Item(136, 11).Value = 18792.666
$ws.Cells.Item(136, 12).Value = 12669.75
$ws.Cells.Item(136, 13).Value = -16242.666
$ws.Cells.Item(136, 14).Value = -17769.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 90910350
$ws.Cells.Item(100, 9).Value = 1180.4
$ws.Cells.Item(100, 11).Value = 2360.8
$ws.Cells.Item(100, 13).Value = -1819.8
$ws.Cells.Item(135, 8).Value = 54609.145
$ws.Cells.Item(135, 10).Value = 54609.145
$ws.Cells.Item(135, 12).Value = 54609.145
$ws.Cells.Item(135, 14).Value = -64749.145
$ws.Cells.Item(136, 8).Value = 19232848
$ws.Cells.Item(136, 9).Value = 23256870
$ws.Cells.Item(136, 10).Value = 6971.5557
$ws.Cells.Item(136, 11).Value = 69770610
$ws.Cells.Item(136, 12).Value = 20914.6671
$ws.Cells.Item(136, 13).Value = -69768060
$ws.Cells.Item(136, 14).Value = -26014.6671
$ws.Cells.Item(141, 8).Value = 69998.08
$ws.Cells.Item(141, 10).Value = 69998.08
$ws.Cells.Item(141, 12).Value = 69998.08
$ws.Cells.Item(141, 14).Value = -80358.08
